$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPGAQ23PU%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T080003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEB0aCXVzLXdlc3QtMiJHMEUCIEmvDbvoJj9lLJNVujVMWKCT%2Ft3cCQd0NxCPBB17JmQ4AiEArnYcU6mg3jz%2F69wqtkm8N6wnA1p00aXVBW%2BccWmjUjgq5gMIhv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCGGI13ACsgw0%2BOByiq6A6Cu6Md3aOAbVqvaqhtiKvr9U90Q8mvsawgmK8ZCnaVflEEHnA3m21IsHTQzH3%2FDgEtXgwCwZSLRCBBx%2FIZBHRYhgidAg0tYv6pKhWTafru%2Fxu%2BSuAdOjaDDTypN5K9%2BHnnNnyyzQs8gRr0aMGRto56zbKKQskQp1GhKIDPaXRwDAd4UNMwT2RcY51wliVZbWpZbjokU306atG6Nhyn7q7MZxPtVFRtDKppFClNXnee4JnlNAofeELCFOqm1Kn%2BOycRf7myl%2Fs6IeeCcYy7SOdJdobT%2BpSWHpVT8kyc1okzt3u%2FL619IuESDsS9EgREQzfQd6W6wLbhJPV09IEche3UICrRRNv%2BI12ncPFfXjDwxzOllwzrEHwzamnkSdCCxDpyolyJNSoSTFQfo%2FKRo4WuBmiPAy8e2ZxaoXbc8%2BrbFhuyWkN%2FRobmAL%2BTnX4ZvTKJstdWA58S82DT2LjZyjg1bcXK49VrkdBB23JWW%2FcHJcne2jfCG1fpNV3EEyKWIaXT6SfRJgFXZ5JdOkOuc9jb78hsSijjmLFG%2FyafCEIBCIJMbFCONpoIkHTxbSoqDPt8UyZoxpnYhsTAwp5SjvwY6pQH3J9Q3mzOu043s6iK%2BxS3RAEMlH9SvLZS3H4s0RA1usyzbXu2Krj8iWXemsURkLZEoVmsf4Nzk9bg8Z60pSPeLpp2xhmfAhlsNA9oIWOSdGXUvvvT8itlLWqOoUb7pajXvjtO1ikX8cTTBzj8KgRZ94xcizgo4C4kbDbcOXoRXB1b72a1MbrSrgF6Azq4r3WejmzAVb2VcFCXL%2Fg5Ob7VyrMz0tto%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2016%3A00%3A03%20GMT&X-Amz-Signature=ff27cbf9d231bea5e16c6db2f91f2dad6944b05468bb49a370182d6882dd9a2f'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPGAQ23PU%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T080003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEB0aCXVzLXdlc3QtMiJHMEUCIEmvDbvoJj9lLJNVujVMWKCT%2Ft3cCQd0NxCPBB17JmQ4AiEArnYcU6mg3jz%2F69wqtkm8N6wnA1p00aXVBW%2BccWmjUjgq5gMIhv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCGGI13ACsgw0%2BOByiq6A6Cu6Md3aOAbVqvaqhtiKvr9U90Q8mvsawgmK8ZCnaVflEEHnA3m21IsHTQzH3%2FDgEtXgwCwZSLRCBBx%2FIZBHRYhgidAg0tYv6pKhWTafru%2Fxu%2BSuAdOjaDDTypN5K9%2BHnnNnyyzQs8gRr0aMGRto56zbKKQskQp1GhKIDPaXRwDAd4UNMwT2RcY51wliVZbWpZbjokU306atG6Nhyn7q7MZxPtVFRtDKppFClNXnee4JnlNAofeELCFOqm1Kn%2BOycRf7myl%2Fs6IeeCcYy7SOdJdobT%2BpSWHpVT8kyc1okzt3u%2FL619IuESDsS9EgREQzfQd6W6wLbhJPV09IEche3UICrRRNv%2BI12ncPFfXjDwxzOllwzrEHwzamnkSdCCxDpyolyJNSoSTFQfo%2FKRo4WuBmiPAy8e2ZxaoXbc8%2BrbFhuyWkN%2FRobmAL%2BTnX4ZvTKJstdWA58S82DT2LjZyjg1bcXK49VrkdBB23JWW%2FcHJcne2jfCG1fpNV3EEyKWIaXT6SfRJgFXZ5JdOkOuc9jb78hsSijjmLFG%2FyafCEIBCIJMbFCONpoIkHTxbSoqDPt8UyZoxpnYhsTAwp5SjvwY6pQH3J9Q3mzOu043s6iK%2BxS3RAEMlH9SvLZS3H4s0RA1usyzbXu2Krj8iWXemsURkLZEoVmsf4Nzk9bg8Z60pSPeLpp2xhmfAhlsNA9oIWOSdGXUvvvT8itlLWqOoUb7pajXvjtO1ikX8cTTBzj8KgRZ94xcizgo4C4kbDbcOXoRXB1b72a1MbrSrgF6Azq4r3WejmzAVb2VcFCXL%2Fg5Ob7VyrMz0tto%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2016%3A00%3A03%20GMT&X-Amz-Signature=ed8cf02f2ce7f77240551501a51dfa5e5575baaf62b145db3160d253a38ec679'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPGAQ23PU%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T080003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEB0aCXVzLXdlc3QtMiJHMEUCIEmvDbvoJj9lLJNVujVMWKCT%2Ft3cCQd0NxCPBB17JmQ4AiEArnYcU6mg3jz%2F69wqtkm8N6wnA1p00aXVBW%2BccWmjUjgq5gMIhv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCGGI13ACsgw0%2BOByiq6A6Cu6Md3aOAbVqvaqhtiKvr9U90Q8mvsawgmK8ZCnaVflEEHnA3m21IsHTQzH3%2FDgEtXgwCwZSLRCBBx%2FIZBHRYhgidAg0tYv6pKhWTafru%2Fxu%2BSuAdOjaDDTypN5K9%2BHnnNnyyzQs8gRr0aMGRto56zbKKQskQp1GhKIDPaXRwDAd4UNMwT2RcY51wliVZbWpZbjokU306atG6Nhyn7q7MZxPtVFRtDKppFClNXnee4JnlNAofeELCFOqm1Kn%2BOycRf7myl%2Fs6IeeCcYy7SOdJdobT%2BpSWHpVT8kyc1okzt3u%2FL619IuESDsS9EgREQzfQd6W6wLbhJPV09IEche3UICrRRNv%2BI12ncPFfXjDwxzOllwzrEHwzamnkSdCCxDpyolyJNSoSTFQfo%2FKRo4WuBmiPAy8e2ZxaoXbc8%2BrbFhuyWkN%2FRobmAL%2BTnX4ZvTKJstdWA58S82DT2LjZyjg1bcXK49VrkdBB23JWW%2FcHJcne2jfCG1fpNV3EEyKWIaXT6SfRJgFXZ5JdOkOuc9jb78hsSijjmLFG%2FyafCEIBCIJMbFCONpoIkHTxbSoqDPt8UyZoxpnYhsTAwp5SjvwY6pQH3J9Q3mzOu043s6iK%2BxS3RAEMlH9SvLZS3H4s0RA1usyzbXu2Krj8iWXemsURkLZEoVmsf4Nzk9bg8Z60pSPeLpp2xhmfAhlsNA9oIWOSdGXUvvvT8itlLWqOoUb7pajXvjtO1ikX8cTTBzj8KgRZ94xcizgo4C4kbDbcOXoRXB1b72a1MbrSrgF6Azq4r3WejmzAVb2VcFCXL%2Fg5Ob7VyrMz0tto%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2016%3A00%3A03%20GMT&X-Amz-Signature=f097885ec62da588233c786e5e2e5fba984db2fbb1b95851828fbe4f73e6e18e'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPGAQ23PU%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T080003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEB0aCXVzLXdlc3QtMiJHMEUCIEmvDbvoJj9lLJNVujVMWKCT%2Ft3cCQd0NxCPBB17JmQ4AiEArnYcU6mg3jz%2F69wqtkm8N6wnA1p00aXVBW%2BccWmjUjgq5gMIhv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCGGI13ACsgw0%2BOByiq6A6Cu6Md3aOAbVqvaqhtiKvr9U90Q8mvsawgmK8ZCnaVflEEHnA3m21IsHTQzH3%2FDgEtXgwCwZSLRCBBx%2FIZBHRYhgidAg0tYv6pKhWTafru%2Fxu%2BSuAdOjaDDTypN5K9%2BHnnNnyyzQs8gRr0aMGRto56zbKKQskQp1GhKIDPaXRwDAd4UNMwT2RcY51wliVZbWpZbjokU306atG6Nhyn7q7MZxPtVFRtDKppFClNXnee4JnlNAofeELCFOqm1Kn%2BOycRf7myl%2Fs6IeeCcYy7SOdJdobT%2BpSWHpVT8kyc1okzt3u%2FL619IuESDsS9EgREQzfQd6W6wLbhJPV09IEche3UICrRRNv%2BI12ncPFfXjDwxzOllwzrEHwzamnkSdCCxDpyolyJNSoSTFQfo%2FKRo4WuBmiPAy8e2ZxaoXbc8%2BrbFhuyWkN%2FRobmAL%2BTnX4ZvTKJstdWA58S82DT2LjZyjg1bcXK49VrkdBB23JWW%2FcHJcne2jfCG1fpNV3EEyKWIaXT6SfRJgFXZ5JdOkOuc9jb78hsSijjmLFG%2FyafCEIBCIJMbFCONpoIkHTxbSoqDPt8UyZoxpnYhsTAwp5SjvwY6pQH3J9Q3mzOu043s6iK%2BxS3RAEMlH9SvLZS3H4s0RA1usyzbXu2Krj8iWXemsURkLZEoVmsf4Nzk9bg8Z60pSPeLpp2xhmfAhlsNA9oIWOSdGXUvvvT8itlLWqOoUb7pajXvjtO1ikX8cTTBzj8KgRZ94xcizgo4C4kbDbcOXoRXB1b72a1MbrSrgF6Azq4r3WejmzAVb2VcFCXL%2Fg5Ob7VyrMz0tto%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2016%3A00%3A03%20GMT&X-Amz-Signature=58e1ed3755e8b53ffe3666e09d0f3fb131d89cdf4e21c72d4ff6c55b04714980'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPGAQ23PU%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T080003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEB0aCXVzLXdlc3QtMiJHMEUCIEmvDbvoJj9lLJNVujVMWKCT%2Ft3cCQd0NxCPBB17JmQ4AiEArnYcU6mg3jz%2F69wqtkm8N6wnA1p00aXVBW%2BccWmjUjgq5gMIhv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCGGI13ACsgw0%2BOByiq6A6Cu6Md3aOAbVqvaqhtiKvr9U90Q8mvsawgmK8ZCnaVflEEHnA3m21IsHTQzH3%2FDgEtXgwCwZSLRCBBx%2FIZBHRYhgidAg0tYv6pKhWTafru%2Fxu%2BSuAdOjaDDTypN5K9%2BHnnNnyyzQs8gRr0aMGRto56zbKKQskQp1GhKIDPaXRwDAd4UNMwT2RcY51wliVZbWpZbjokU306atG6Nhyn7q7MZxPtVFRtDKppFClNXnee4JnlNAofeELCFOqm1Kn%2BOycRf7myl%2Fs6IeeCcYy7SOdJdobT%2BpSWHpVT8kyc1okzt3u%2FL619IuESDsS9EgREQzfQd6W6wLbhJPV09IEche3UICrRRNv%2BI12ncPFfXjDwxzOllwzrEHwzamnkSdCCxDpyolyJNSoSTFQfo%2FKRo4WuBmiPAy8e2ZxaoXbc8%2BrbFhuyWkN%2FRobmAL%2BTnX4ZvTKJstdWA58S82DT2LjZyjg1bcXK49VrkdBB23JWW%2FcHJcne2jfCG1fpNV3EEyKWIaXT6SfRJgFXZ5JdOkOuc9jb78hsSijjmLFG%2FyafCEIBCIJMbFCONpoIkHTxbSoqDPt8UyZoxpnYhsTAwp5SjvwY6pQH3J9Q3mzOu043s6iK%2BxS3RAEMlH9SvLZS3H4s0RA1usyzbXu2Krj8iWXemsURkLZEoVmsf4Nzk9bg8Z60pSPeLpp2xhmfAhlsNA9oIWOSdGXUvvvT8itlLWqOoUb7pajXvjtO1ikX8cTTBzj8KgRZ94xcizgo4C4kbDbcOXoRXB1b72a1MbrSrgF6Azq4r3WejmzAVb2VcFCXL%2Fg5Ob7VyrMz0tto%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2016%3A00%3A03%20GMT&X-Amz-Signature=f68225b21a4facf268086e7575c825d2e57ff73aa20329a5082fc094f7428a73'
